$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH126"
$ws.Range("C2").Value = "ROBBEN ISLANDERS TAKE A LOOK AT THE POLITICAL SITUATION IN NATAL"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

$dataRange = $ws.Range("A2:H2")
$dataRange.Font.Name = "Calibri"
$dataRange.Font.Size = 10
